# Apply cryptocurrency price/volume updates (cryptos.xlsx) as of Mon May 15 16:39:31 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Price column cells whose new values look like plain decimal numbers
# as Text, so Excel keeps them as strings (matching the source data) instead of
# silently converting them to floating point numbers.
$textFixCells = @("D4","D5","D7","D8","D9","D10","D11","D13","D14","D15","D16","D18","D19","D20","D22","D23","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D38","D39","D40","D42","D43","D44","D45","D46","D47","D48","D49","D51")
foreach ($addr in $textFixCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update cell values row by row to match the refreshed crypto data feed
$ws.Range("D2").Value = "27.722.11"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "1.853.11"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  -2.07%  "
$ws.Range("D5").Value = "319.10"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("E6").Value = "  -1.85%  "
$ws.Range("D7").Value = "0.4308"
$ws.Range("E7").Value = "  -2.21%  "
$ws.Range("D8").Value = "0.3748"
$ws.Range("E8").Value = "  -1.77%  "
$ws.Range("D9").Value = "0.07339"
$ws.Range("E9").Value = "  -1.48%  "
$ws.Range("D10").Value = "0.8779"
$ws.Range("E10").Value = "  -1.30%  "
$ws.Range("D11").Value = "21.60"
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").Value = "1.847.62"
$ws.Range("E12").Value = "  -1.11%  "
$ws.Range("D13").Value = "6.727"
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("D14").Value = "5.437"
$ws.Range("E14").Value = "  -2.29%  "
$ws.Range("D15").Value = "0.07103"
$ws.Range("E15").Value = "  -1.21%  "
$ws.Range("D16").Value = "88.97"
$ws.Range("E16").Value = "  +3.89%  "
$ws.Range("E17").Value = "  -1.99%  "
$ws.Range("D18").Value = "0.000008992"
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("D19").Value = "1.009"
$ws.Range("E19").Value = "  -1.80%  "
$ws.Range("D20").Value = "15.48"
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("D21").Value = "27.742.45"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").Value = "5.215"
$ws.Range("E22").Value = "  -1.94%  "
$ws.Range("D23").Value = "11.07"
$ws.Range("E23").Value = "  -1.93%  "
$ws.Range("D24").Value = "2.079.92"
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("D25").Value = "1.986"
$ws.Range("E25").Value = "  -2.38%  "
$ws.Range("D26").Value = "155.40"
$ws.Range("E26").Value = "  -2.01%  "
$ws.Range("D27").Value = "18.63"
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("D28").Value = "2.189"
$ws.Range("E28").Value = "  +8.98%  "
$ws.Range("D29").Value = "5.369"
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("D30").Value = "118.98"
$ws.Range("E30").Value = "  +0.84%  "
$ws.Range("D31").Value = "0.08940"
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("D32").Value = "1.231"
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("D33").Value = "0.7779"
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("D34").Value = "4.551"
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("E35").Value = "  -3.40%  "
$ws.Range("E36").Value = "  -1.90%  "
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01983"
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.05352"
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").Value = "7.344"
$ws.Range("E40").Value = "  +6.01%  "
$ws.Range("E41").Value = "  +1.66%  "
$ws.Range("D42").Value = "0.1692"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").Value = "0.5135"
$ws.Range("E43").Value = "  -1.92%  "
$ws.Range("D44").Value = "8.812"
$ws.Range("E44").Value = "  -1.16%  "
$ws.Range("D45").Value = "10.74"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").Value = "108.12"
$ws.Range("E46").Value = "  -2.66%  "
$ws.Range("D47").Value = "0.4789"
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("D48").Value = "0.06474"
$ws.Range("E48").Value = "  -2.00%  "
$ws.Range("D49").Value = "1.691"
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("E50").Value = "  -2.02%  "
$ws.Range("D51").Value = "1.846"
$ws.Range("E51").Value = "  -4.30%  "
